# cryptos.xlsx refresh - updates coin price/volume(1h) snapshot values
# (and a couple of rank-swapped rows) pulled from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "90.008.15"
$ws.Range("E2").Value = "  +2.39%  "
# Row 3
$ws.Range("D3").Value = "3.184.97"
$ws.Range("E3").Value = "  -2.27%  "
# Row 4
$ws.Range("E4").Value = "  -0.29%  "
# Row 5
$ws.Range("D5").Formula = "'214.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.46%  "
# Row 6
$ws.Range("D6").Formula = "'618.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.43%  "
# Row 7
$ws.Range("D7").Formula = "'0.395"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.18%  "
# Row 8
$ws.Range("E8").Value = "  -3.60%  "
# Row 9
$ws.Range("D9").Formula = "'0.999"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.13%  "
# Row 10
$ws.Range("D10").Value = "3.178.81"
$ws.Range("E10").Value = "  -2.74%  "
# Row 11
$ws.Range("D11").Formula = "'0.577"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.09%  "
# Row 12
$ws.Range("D12").Formula = "'0.177"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.70%  "
# Row 13
$ws.Range("D13").Formula = "'0.0000256"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.38%  "
# Row 14
$ws.Range("D14").Value = "89.713.85"
$ws.Range("E14").Value = "  +2.25%  "
# Row 15
$ws.Range("D15").Value = "3.767.79"
$ws.Range("E15").Value = "  -2.65%  "
# Row 16
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Formula = "'32.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.05%  "
# Row 17
$ws.Range("B17").Value = "Toncoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D17").Formula = "'5.26"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.40%  "
# Row 18
$ws.Range("D18").Value = "3.189.24"
$ws.Range("E18").Value = "  -2.49%  "
# Row 19
$ws.Range("D19").Formula = "'3.29"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.04%  "
# Row 20
$ws.Range("B20").Value = "PEPE"
$ws.Range("C20").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D20").Formula = "'0.0000200"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +40.85%  "
# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Formula = "'13.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.63%  "
# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Formula = "'437.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.17%  "
# Row 23
$ws.Range("D23").Formula = "'8.62"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.05%  "
# Row 24
$ws.Range("D24").Formula = "'5.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.35%  "
# Row 25
$ws.Range("D25").Formula = "'5.15"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.45%  "
# Row 26
$ws.Range("D26").Formula = "'11.67"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.26%  "
# Row 27
$ws.Range("D27").Value = "3.341.99"
$ws.Range("E27").Value = "  -3.17%  "
# Row 28
$ws.Range("D28").Formula = "'75.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.77%  "
# Row 29
$ws.Range("D29").Formula = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.10%  "
# Row 30
$ws.Range("D30").Formula = "'0.167"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.86%  "
# Row 31
$ws.Range("D31").Formula = "'0.999"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.28%  "
# Row 32
$ws.Range("D32").Formula = "'4.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +25.26%  "
# Row 33
$ws.Range("D33").Formula = "'8.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.92%  "
# Row 34
$ws.Range("D34").Formula = "'536.50"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.86%  "
# Row 35
$ws.Range("D35").Formula = "'7.05"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.82%  "
# Row 36
$ws.Range("D36").Formula = "'1.86"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.33%  "
# Row 37
$ws.Range("D37").Formula = "'1.27"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -8.48%  "
# Row 38
$ws.Range("D38").Formula = "'22.05"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.09%  "
# Row 39
$ws.Range("E39").Value = "  +2.23%  "
# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Formula = "'0.127"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.88%  "
# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Formula = "'0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.20%  "
# Row 42
$ws.Range("E42").Value = "  +0.07%  "
# Row 43
$ws.Range("D43").Formula = "'1.94"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.19%  "
# Row 44
$ws.Range("D44").Formula = "'0.374"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.27%  "
# Row 45
$ws.Range("D45").Formula = "'148.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.05%  "
# Row 46
$ws.Range("D46").Formula = "'43.67"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.47%  "
# Row 47
$ws.Range("D47").Formula = "'172.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.24%  "
# Row 48
$ws.Range("D48").Formula = "'0.124"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -10.65%  "
# Row 49
$ws.Range("D49").Formula = "'1.24"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.50%  "
# Row 50
$ws.Range("D50").Formula = "'4.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.93%  "
# Row 51
$ws.Range("D51").Formula = "'0.609"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.01%  "

Write-Host "Applied cryptos.xlsx updates"
